# Auto-generated edit script applying the crypto price-table update
# described in the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.856.89'
$ws.Range("E2").Value = '  +2.42%  '
$ws.Range("D3").Value = '2.948.18'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.62'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.43'
$ws.Range("E6").Value = '  +2.54%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '2.947.35'
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.06'
$ws.Range("E10").Value = '  +1.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.149'
$ws.Range("E11").Value = '  +5.12%  '
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000234'
$ws.Range("E13").Value = '  +4.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.53'
$ws.Range("E14").Value = '  -2.19%  '
$ws.Range("D16").Value = '3.433.80'
$ws.Range("D17").Value = '62.839.49'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.69'
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").Value = '2.945.11'
$ws.Range("E19").Value = '  +0.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '438.22'
$ws.Range("E20").Value = '  +1.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.40'
$ws.Range("E21").Value = '  -1.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.665'
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.28'
$ws.Range("E24").Value = '  +3.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.68'
$ws.Range("E25").Value = '  -0.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.13'
$ws.Range("E26").Value = '  -2.21%  '
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  +1.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.32'
$ws.Range("E30").Value = '  +6.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.61'
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("D32").Value = '0.0₃0989'
$ws.Range("E32").Value = '  +13.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.36'
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.991'
$ws.Range("E36").Value = '  -1.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.61'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.59'
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.02'
$ws.Range("E40").Value = '  +1.30%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.118'
$ws.Range("E41").Value = '  -3.90%  '
$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.44'
$ws.Range("E42").Value = '  -0.70%  '
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.54'
$ws.Range("E44").Value = '  -5.99%  '
$ws.Range("D45").Value = '2.690.57'
$ws.Range("E45").Value = '  -0.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '134.99'
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0338'
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '357.70'
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.67'
$ws.Range("E51").Value = '  -3.78%  '
